$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K shifts to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

$rowData = @(
    @{R=7; D="43465"; E="43373"},
    @{R=8; D="279400"; E="269400"},
    @{R=9; D="207400"; E="200400"},
    @{R=10; D="72000"; E="69000"},
    @{R=11; D="NULL"; E="NULL"},
    @{R=12; D="NA"; E="NA"},
    @{R=13; D="0"; E="0"},
    @{R=14; D="200"; E="400"},
    @{R=15; D="0"; E="0"},
    @{R=16; D="NULL"; E="NULL"},
    @{R=17; D="246000"; E="239000"},
    @{R=18; D="33400"; E="30400"},
    @{R=19; D="NULL"; E="NULL"},
    @{R=20; D="-100"; E="0"},
    @{R=21; D="42700"; E="39800"},
    @{R=22; D="2100"; E="2200"},
    @{R=23; D="31200"; E="28200"},
    @{R=24; D="-1000"; E="6500"},
    @{R=25; D="0"; E="0"},
    @{R=26; D="32200"; E="21700"},
    @{R=27; D="32200"; E="21700"},
    @{R=28; D="0"; E="0"},
    @{R=29; D="300"; E="NA"},
    @{R=30; D="0"; E="0"},
    @{R=31; D="0"; E="0"},
    @{R=32; D="100"; E="0"},
    @{R=33; D="32500"; E="21700"},
    @{R=34; D="0"; E="0"},
    @{R=35; D="32500"; E="21700"},
    @{R=38; D="43465"; E="43373"},
    @{R=39; D="NULL"; E="NULL"},
    @{R=40; D="NULL"; E="NULL"},
    @{R=41; D="37400"; E="35000"},
    @{R=42; D="0"; E="0"},
    @{R=43; D="124400"; E="117800"},
    @{R=44; D="157300"; E="153600"},
    @{R=45; D="9400"; E="8000"},
    @{R=46; D="328500"; E="314400"},
    @{R=47; D="0"; E="0"},
    @{R=48; D="158600"; E="159900"},
    @{R=49; D="518200"; E="521400"},
    @{R=50; D="0"; E="0"},
    @{R=51; D="0"; E="0"},
    @{R=52; D="18500"; E="13600"},
    @{R=53; D="0"; E="0"},
    @{R=54; D="1023800"; E="1009300"},
    @{R=55; D="NULL"; E="NULL"},
    @{R=56; D="NULL"; E="NULL"},
    @{R=57; D="66100"; E="66000"},
    @{R=58; D="200"; E="300"},
    @{R=59; D="92500"; E="88800"},
    @{R=60; D="158800"; E="155100"},
    @{R=61; D="209900"; E="221900"},
    @{R=62; D="125000"; E="120500"},
    @{R=63; D="0"; E="0"},
    @{R=64; D="0"; E="0"},
    @{R=65; D="0"; E="0"},
    @{R=66; D="493700"; E="497500"},
    @{R=67; D="NULL"; E="NULL"},
    @{R=68; D="0"; E="0"},
    @{R=69; D="0"; E="0"},
    @{R=70; D="0"; E="0"},
    @{R=71; D="0"; E="0"},
    @{R=72; D="432500"; E="394300"},
    @{R=73; D="0"; E="0"},
    @{R=74; D="0"; E="0"},
    @{R=75; D="0"; E="0"},
    @{R=76; D="530100"; E="511800"},
    @{R=77; D="0"; E="0"},
    @{R=80; D="43465"; E="43373"},
    @{R=81; D="32500"; E="21700"},
    @{R=82; D="NULL"; E="NULL"},
    @{R=83; D="9400"; E="9400"},
    @{R=84; D="0"; E="0"},
    @{R=85; D="0"; E="0"},
    @{R=86; D="0"; E="0"},
    @{R=87; D="0"; E="0"},
    @{R=88; D="0"; E="0"},
    @{R=89; D="21000"; E="34000"},
    @{R=90; D="NULL"; E="NULL"},
    @{R=91; D="-4000"; E="-3100"},
    @{R=92; D="0"; E="0"},
    @{R=93; D="0"; E="0"},
    @{R=94; D="-3900"; E="-3200"},
    @{R=95; D="NULL"; E="NULL"},
    @{R=96; D="-4900"; E="-4800"},
    @{R=97; D="0"; E="0"},
    @{R=98; D="0"; E="0"},
    @{R=99; D="0"; E="0"},
    @{R=100; D="-14600"; E="-31700"},
    @{R=101; D="-100"; E="-100"},
    @{R=102; D="2400"; E="-1000"},
)

foreach ($item in $rowData) {
    $r = $item.R
    # Copy number format / style from column F (the shifted former column D) onto D and E
    $ws.Range("F$r").Copy()
    $ws.Range("D${r}:E$r").PasteSpecial(-4122)

    if ($item.D -ne "NULL") {
        $ws.Range("D$r").Value = $item.D
    }
    if ($item.E -ne "NULL") {
        $ws.Range("E$r").Value = $item.E
    }
}

$excel.CutCopyMode = 0